# Update Name of Algo
# Apply corrected values to the "result_data_KNN" sheet, matching the
# updated algorithm output (column A and column C values changed for a
# number of rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.228
$ws.Range("A12").Value = -21.721
$ws.Range("C14").Value = -12.139
$ws.Range("C26").Value = -12.805
$ws.Range("C31").Value = -13.27
$ws.Range("A32").Value = -21.757
$ws.Range("C35").Value = -12.762
$ws.Range("A36").Value = -20.178
$ws.Range("C37").Value = -13.341
$ws.Range("A38").Value = -19.741
$ws.Range("C45").Value = -12.898
$ws.Range("A46").Value = -21.924
$ws.Range("A54").Value = -22.155
$ws.Range("A55").Value = -22.219
$ws.Range("C57").Value = -13.829
$ws.Range("A67").Value = -21.495
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("A91").Value = -21.522
$ws.Range("A99").Value = -20.428
$ws.Range("C100").Value = -13.018
$ws.Range("C102").Value = -12.87
